$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column J
$ws.Range("J1").Value = "xheight_to_size"

# Update word/nonword trial counts on row 2
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 20

# Add new xheight_to_size values for rows 2-7
$ws.Range("J2").Value = 2.0699999999999998
$ws.Range("J3").Value = 1.9824999999999999
$ws.Range("J4").Value = 1.96
$ws.Range("J5").Value = 2.0699999999999998
$ws.Range("J6").Value = 1.9824999999999999
$ws.Range("J7").Value = 1.96

# Update the selected cell to match the target view
$ws.Range("E4").Select()
